$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.881.19"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.272.05"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D5").Value = "303.94"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "93.27"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").Value = "0.531"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "32.77"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "53.57"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "0.0797"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "6.71"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "2.623.94"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "14.33"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "2.277.30"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "0.780"
$ws.Range("E18").Value = "  +3.09%  "
$ws.Range("D19").Value = "41.799.83"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "12.89"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "67.28"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "244.31"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "24.10"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "9.55"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  -5.16%  "
$ws.Range("D31").Value = "35.22"
$ws.Range("E31").Value = "  +3.03%  "
$ws.Range("D32").Value = "160.54"
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "0.0745"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").Value = "3.04"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.106"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "16.94"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("D39").Value = "2.37"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "1.81"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").Value = "3.95"
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "19.98"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.017.62"
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").Value = "10.43"
$ws.Range("E46").Value = "  +3.30%  "
$ws.Range("D47").Value = "2.14"
$ws.Range("E47").Value = "  +7.75%  "
$ws.Range("D48").Value = "2.91"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").Value = "53.52"
$ws.Range("E49").Value = "  +3.72%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "72.60"
$ws.Range("E51").Value = "  +3.10%  "
